$wb = $excel.ActiveWorkbook

# Update "Correspond Handoff Datetime" (E3) and "Correspond Handback DateTime" (H3)
# for the "1c27b64b-..." report row on both locale report sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-13 18:48:07"
$wsZhCn.Range("H3").Value = "2016-03-13 18:48:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-13 18:48:10"
$wsDeDe.Range("H3").Value = "2016-03-13 18:48:32"
